$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "example"

# Merge name.first + name.last into a single "name" column (column B),
# shifting address/facebook/gender columns left by one.
$ws.Range("B1").Value = "name"
$ws.Range("B2").Value = "Jason Humphrey"
$ws.Range("B3").Value = "Tyler Sorber"

$ws.Range("C1").Value = "address,1_location"
$ws.Range("C2").Value = "texas"
$ws.Range("C3").Value = "ohio"

$ws.Range("D1").Value = "address,2_location"
$ws.Range("D2").Value = "florida"
$ws.Range("D3").Value = "california"

$ws.Range("E1").Value = "facebook.id"
$ws.Range("E2").Value = 1232132121
$ws.Range("E3").Value = 1232342143

$ws.Range("F1").Value = "facebook.imageurl"
$ws.Range("F2").Value = "https://graph.facebook.com/1232132121/picture?height=200&width=200&type=square"
$ws.Range("F3").Value = "https://graph.facebook.com/1232342143/picture?height=200&width=200&type=square"

$ws.Range("G1").Value = "facebook.image"
$ws.Range("G2").Value = "Yes"
$ws.Range("G3").Value = "Yes"

$ws.Range("H1").Value = "gender"
$ws.Range("H2").Value = "Male"
$ws.Range("H3").Value = "Male"

# Remove now-unused column I (was facebook.image / gender before the shift)
$ws.Range("I1:I3").Clear()

# Update the selection to match the new active cell
$ws.Range("H5").Select()
